$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Hunk 1: merge the two runs in the first paragraph ("...employment ends"
# + ".") into a single run carrying the trailing period, and drop the
# _GoBack bookmark that used to sit at the end of that paragraph (it will
# be re-created further down, next to its new home).
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("employment ends.", $false, $false, $false, $false, $false, $true, 1, $false, "employment ends.", 2) | Out-Null

$d.Bookmarks("_GoBack").Delete()

# ---------------------------------------------------------------------
# Hunk 2: "Document of Record for a policy change request" becomes
# "Document of Record for" + " a workforce member termination." (kept as
# two separate runs with identical formatting, mirroring a real edit),
# with a fresh, collapsed _GoBack bookmark right after the new text.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Document of Record for a policy change request", $false, $false, $false, $false, $false, $true, 1, $false, "Document of Record for a workforce member termination.", 2) | Out-Null

$splitRng = $d.Content
$splitRng.Find.Execute(" a workforce member termination.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
# Force a run boundary right before " a workforce member termination." by
# toggling a formatting property on just that sub-range, then restoring it.
$splitRng.Bold = 1
$splitRng.Bold = 0

$splitRng.Collapse(0)
$d.Bookmarks.Add("_GoBack", $splitRng) | Out-Null

# ---------------------------------------------------------------------
# Hunk 3: "   Voluntary      " + "or      " (with proofErr gramStart/
# gramEnd wrapping the second run) collapse into a single run.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("   Voluntary      or      ", $false, $false, $false, $false, $false, $true, 1, $false, "   Voluntary      or      ", 2) | Out-Null

# ---------------------------------------------------------------------
# Hunk 4: the Google Apps Access bullet collapses back into one run,
# removing the spellcheck proofErr wrappers around "gDrive" / "etc".
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Terminate Employee" + [char]0x2019 + "s Google Apps Access (email, calendar, gDrive, etc).", $false, $false, $false, $false, $false, $true, 1, $false, "Terminate Employee" + [char]0x2019 + "s Google Apps Access (email, calendar, gDrive, etc).", 2) | Out-Null

# ---------------------------------------------------------------------
# Hunk 5: the Github Access bullet collapses back into one run, removing
# the spellcheck proofErr wrappers around "Github".
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Terminate Employee" + [char]0x2019 + "s Github Access", $false, $false, $false, $false, $false, $true, 1, $false, "Terminate Employee" + [char]0x2019 + "s Github Access", 2) | Out-Null
